$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value (7 -> 10)
$ws.Range("B2").Value = 10

# Update row 3 (A3: 2 -> 1, B3: 4 -> 2)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2

# Update row 4 (A4: 3 -> 2, B4 stays 1)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1

# Remove old row 5 entirely (A5=1, B5=1), shifting rows up so the range becomes A1:B4
$ws.Rows("5:5").Delete()
